# Main.xlsx / Feuille1 (BOM) update
#  - swap the part number + description of the MLCC capacitor row (row 17)
#  - filter the "Ord. Qty" column (B) down to 1 and 4, which hides every
#    other data row
#  - leave the resulting selection on A2:B23 (matches the post-edit file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM edit: replaced capacitor part -------------------------------
$ws.Range("A17").Value = "FG18C0G1H331JNT00"
$ws.Range("C17").Value = "Multilayer Ceramic Capacitors MLCC - Leaded RAD 50V 330pF C0G 5% LS:2.5mm"

# --- AutoFilter on column B (Ord. Qty), keep rows where B is 1 or 4 ---
$ws.Range("A1:E23").AutoFilter(2, @("1", "4"), 7)

# Excel always writes a hidden workbook-level _FilterDatabase name once a
# filter is applied to a range - recreate it so the saved file matches.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Feuille1!`$A`$1:`$E`$23")
$filterName.Visible = $false

# --- Selection left on A2:B23 after the edit --------------------------
$ws.Range("A2:B23").Select()
